# Append 5 new rows (16-20) of "decision_pending" data to the
# "decision_pending" sheet, grow the Table1 ListObject to cover them,
# and leave the active cell on A18 (matching the author's final selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("decision_pending")

$data = @(
    ,@("Medical Subject", "Medical Subject: Topic", "yes", "not indexed", "n", "subject_heading_mesh", ".", ".", ".", ".", ".", "map", "https://trlnmain.atlassian.net/browse/TD-131")
    ,@("Indexed Notes", "Indexed Notes", "no", "Keyword", "y", "tbd", "y", ".", ".", ".", "split up into other fields", "map", "https://trlnmain.atlassian.net/browse/TD-160")
    ,@("Notes", "Notes", "no", "not indexed", "y", "tbd", "y", ".", ".", ".", "split up into other fields", "map", "https://trlnmain.atlassian.net/browse/TD-160")
    ,@("Collective Uniform Title", "Collective Uniform Title", "no", "not indexed", "y", "work_name_title", "y", "Related items > ??", ".", "not indexed", "Title portion indexed in indexed-only title_uniform.  Author portion indexed in author_main. Stored for display and creation of hyperlink", "map", "https://trlnmain.atlassian.net/browse/TD-442")
    ,@("Main Uniform Title", "Main Uniform Title", "no", "Keyword, Keyword_Brief, Title, Title_Brief", "y", "work_name_title", "y", "Related works > Other versions of this work", ".", "title; author; keyword", "Title portion indexed in indexed-only title_uniform.  Author portion indexed in author_main. Stored for display and creation of hyperlink", "map", "https://trlnmain.atlassian.net/browse/TD-392")
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $rowValues[$j]
    }
}

# Grow the worksheet's table (Table1) so it covers the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:M20"))

# Match the author's final selection state.
$ws.Range("A18").Select()
